$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.147.80"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "2.484.84"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").Value = "2.483.75"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.332"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").Value = "67.032.72"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").Value = "2.497.93"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("B24").Value = "NEARProtocol"
$ws.Range("C24").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.87%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "2.607.92"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").Value = "0.0₃0901"
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "512.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.73%  "
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("E34").Value = "  -3.83%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.118"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("E40").Value = "  -5.96%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  -3.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.328"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.515"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.34%  "
$ws.Range("D50").Value = "0.0₆0251"
$ws.Range("E50").Value = "  -6.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0728"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.19%  "
